# Update project plan and estimates
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Estimates")

# --- Mark existing stories as Completed (adds Completed checkbox + recalculates
#     Completed Points / Completed Hours via the table's calculated columns) ---
$completedRows = @(74,110,111,112,113,114,120,124,125,126,127)
foreach ($r in $completedRows) {
    $ws.Cells.Item($r,5).Value = $true
}

# --- Fix the typo in the story name for row 126 ---
$ws.Cells.Item(126,2).Value = "Re-add compilerconfig.json"

# --- Insert 12 new rows into the table (pushes the totals/summary block down) ---
$ws.Rows("131:142").Insert()

$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:G142"))

# --- New story rows ---
$newRows = @(
    @{ Row=131; Epic="Tech Debt";        Story="Add Notification Code to Sorry Page";    Points=1; Completed=$true  },
    @{ Row=132; Epic="Manage Locations"; Story="Metro Areas";                             Points=1; Completed=$true  },
    @{ Row=133; Epic="Locations";        Story="Polaris Donations Page";                  Points=3; Completed=$false },
    @{ Row=134; Epic="New Stories";      Story="Send confirmation email for Bed Request"; Points=1; Completed=$false },
    @{ Row=135; Epic="Tech Debt";        Story="Make seeding environment specific";       Points=1; Completed=$false },
    @{ Row=136; Epic="Tech Debt";        Story="Import Bed Requests";                     Points=3; Completed=$false },
    @{ Row=137; Epic="Tech Debt";        Story="Add Database Indexes";                    Points=1; Completed=$false },
    @{ Row=138; Epic="Tech Debt";        Story="Verify Mobile for all pages";              Points=1; Completed=$false },
    @{ Row=139; Epic="Tech Debt";        Story="Refactor Grid Persistence";                Points=2; Completed=$false },
    @{ Row=140; Epic="Tech Debt";        Story="Add try catch and alert";                  Points=3; Completed=$false },
    @{ Row=141; Epic="National Pages";   Story="National History of Bed Brigade Page";     Points=1; Completed=$false },
    @{ Row=142; Epic="National Pages";   Story="National Donations Page";                  Points=1; Completed=$false }
)

foreach ($item in $newRows) {
    $r = $item.Row
    $ws.Cells.Item($r,1).Value = $item.Epic
    $ws.Cells.Item($r,2).Value = $item.Story
    $ws.Cells.Item($r,3).Value = $item.Points
    $ws.Cells.Item($r,4).Formula = "=VLOOKUP(C$r,Points!`$A`$1:`$C`$6,3,FALSE)"
    if ($item.Completed) {
        $ws.Cells.Item($r,5).Value = $true
    }
    $ws.Cells.Item($r,6).Formula = "=IF(Table1[[#This Row],[Completed]],Table1[[#This Row],[Points]],0)"
    $ws.Cells.Item($r,7).Formula = "=IF(Table1[[#This Row],[Completed]],Table1[[#This Row],[Estimated Hours]],0)"
}

# Match the styling Excel applies to newly-appended table rows: the Estimated
# Hours column centered, all three formula columns with an explicit General
# number format.
$ws.Range("D131:D142").HorizontalAlignment = -4108
$ws.Range("D131:D142").NumberFormat = "General"
$ws.Range("F131:G142").NumberFormat = "General"

# --- Move selection / scroll position to match the saved view ---
$ws.Activate()
$ws.Range("D148").Select()
